$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 5294
$ws1.Range("F3").Value = 573
$ws1.Range("F4").Value = 10803
$ws1.Range("F6").Value = 576
$ws1.Range("F7").Value = 151
$ws1.Range("F8").Value = 196
$ws1.Range("F9").Value = 914

# Sheet "演出" (Performances)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 15
$ws2.Range("F4").Value = 21
$ws2.Range("F6").Value = 4

# Sheet "全部类型" (All types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 15
$ws4.Range("F4").Value = 5294
$ws4.Range("F5").Value = 573
$ws4.Range("F6").Value = 21
$ws4.Range("F7").Value = 10803
$ws4.Range("F9").Value = 576
$ws4.Range("F10").Value = 151
$ws4.Range("F12").Value = 4
$ws4.Range("F13").Value = 196
$ws4.Range("F14").Value = 914
